$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
# Overview sheet: zh-cn / de-de status cells
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: Status column
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de sheet: Status column
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the Status columns (shorter text no longer needs the old width) ---
$newWidth = 12.5

$wsOverview.Range("E:E").ColumnWidth = $newWidth
$wsOverview.Range("F:F").ColumnWidth = $newWidth

$wsZhCn.Range("C:C").ColumnWidth = $newWidth

$wsDeDe.Range("C:C").ColumnWidth = $newWidth
